$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Next.js, Django, Tailwind CSS" -> split into three runs:
#   "Next.js" | "," | " Django, Tailwind CSS"
# The run formatting (Merriweather, not-bold, size 16/18) must stay identical
# on all three pieces, so we force Word to break the run at the two
# boundaries by toggling Font.Bold off/on back to its original value - this
# causes the OOXML serializer to split the underlying <w:r> without actually
# changing any visible formatting.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Next.js, Django, Tailwind CSS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $rng1.Start

# boundary after "Next.js" (7 chars) -> splits "Next.js" | ",Django, Tailwind CSS"
$b1a = $d.Range($start1 + 7, $start1 + 7 + 1)
$b1a.Font.Bold = 1
$b1a.Font.Bold = 0

# boundary after "," (1 more char) -> splits "," | " Django, Tailwind CSS"
$b1b = $d.Range($start1 + 7 + 1, $start1 + 29)
$b1b.Font.Bold = 1
$b1b.Font.Bold = 0

# ---------------------------------------------------------------------------
# Change 2: "Zappa, Beautiful Soup, Selenium" -> split into three runs and
# insert the new " React," piece:
#   "Zappa," | " React," | " Beautiful Soup, Selenium"
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Zappa, Beautiful Soup, Selenium", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $rng2.Start

# Insert " React," right after "Zappa," (6 chars in)
$insPt = $d.Range($start2 + 6, $start2 + 6)
$insPt.InsertAfter(" React,")

# Now force the split boundaries so we end up with three separate runs:
#   [start2, start2+6)      -> "Zappa,"
#   [start2+6, start2+13)   -> " React,"
#   [start2+13, start2+39)  -> " Beautiful Soup, Selenium"
$b2a = $d.Range($start2, $start2 + 6)
$b2a.Font.Bold = 1
$b2a.Font.Bold = 0

$b2b = $d.Range($start2 + 6, $start2 + 13)
$b2b.Font.Bold = 1
$b2b.Font.Bold = 0

Write-Host "Done"
